$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date field text on every
#    slide layout footer placeholder (31/08/2021 -> 10/11/2021).
# ------------------------------------------------------------------
$oldDate = "31/08/2021"
$newDate = "10/11/2021"

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shape = $layout.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# ------------------------------------------------------------------
# 2) Remove the leftover "API dos correios" picture + rectangle
#    callout from the single content slide.
# ------------------------------------------------------------------
$s = $p.Slides.Item(1)
for ($k = $s.Shapes.Count; $k -ge 1; $k--) {
    $shape = $s.Shapes.Item($k)
    $isTargetPic = ($shape.Type -eq 13) -and ($shape.Name -eq "Picture 18") -and ($shape.Id -eq 52)
    $isTargetRect = ($shape.Name -eq "Retângulo 74") -and ($shape.Id -eq 75)
    if ($isTargetPic -or $isTargetRect) {
        $shape.Delete()
    }
}
